$d = $word.ActiveDocument

# Update the date/day heading
$d.Content.Find.Execute("2023-06-07 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-08 Thursday", 2) | Out-Null

# New values for each of the 100 table cells, in row-major order (row 1 col 1..5, row 2 col 1..5, ...)
$newValues = @(
    "70-65=",
    "39-21=",
    "27+43=",
    "2-2=",
    "66+14=",
    "94-80=",
    "79-48=",
    "52+19=",
    "51+4=",
    "84-6=",
    "75-34=",
    "10+18=",
    "86-7=",
    "59-43=",
    "78-21=",
    "93-79=",
    "22+72=",
    "2+85=",
    "76-21=",
    "67+30=",
    "49+46=",
    "4+5=",
    "70+23=",
    "72+14=",
    "11+43=",
    "45-42=",
    "69+2=",
    "93-45=",
    "46+8=",
    "6+53=",
    "7+72=",
    "25-12=",
    "69-54=",
    "16+26=",
    "86-8=",
    "6+59=",
    "86-80=",
    "7+72=",
    "82+9=",
    "22-12=",
    "82-78=",
    "93-10=",
    "53-3=",
    "77-72=",
    "63-63=",
    "82-37=",
    "70-42=",
    "90-48=",
    "7+38=",
    "88-27=",
    "56+22=",
    "85-58=",
    "52-44=",
    "77+2=",
    "22-17=",
    "94-39=",
    "69+4=",
    "25-23=",
    "61-21=",
    "68-32=",
    "47-40=",
    "21+6=",
    "43+29=",
    "4+33=",
    "74-16=",
    "21+66=",
    "65+30=",
    "26-12=",
    "80-20=",
    "79-0=",
    "57+39=",
    "58+5=",
    "3-3=",
    "62+6=",
    "30+4=",
    "77-53=",
    "2+42=",
    "60-34=",
    "20+13=",
    "88-8=",
    "95-90=",
    "17-3=",
    "98-2=",
    "59+24=",
    "11+33=",
    "54-15=",
    "74-68=",
    "24+44=",
    "86-44=",
    "74-4=",
    "23+41=",
    "28+22=",
    "97-87=",
    "32+21=",
    "43-32=",
    "17+25=",
    "56-54=",
    "11+22=",
    "79+8=",
    "20+29="
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Output "Updated $idx cells"
